$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "44.219.40"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.244.37"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "306.42"
$ws.Range("E5").Value = "  -2.61%  "
$ws.Range("D6").Value = "95.06"
$ws.Range("E6").Value = "  -3.53%  "
$ws.Range("D7").Value = "0.573"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "34.78"
$ws.Range("E10").Value = "  -3.16%  "
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "7.20"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "2.586.80"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "2.328.56"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("D16").Value = "0.832"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "13.55"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "43.960.11"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "0.0₃0965"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "6.40"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "12.12"
$ws.Range("E21").Value = "  -7.40%  "
$ws.Range("D22").Value = "65.61"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "239.04"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "2.01"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "9.94"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "38.35"
$ws.Range("E28").Value = "  +5.51%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("D30").Value = "20.07"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "5.88"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").Value = "153.35"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "0.0796"
$ws.Range("E33").Value = "  -4.00%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "3.20"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("E36").Value = "  +2.26%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -6.93%  "
$ws.Range("D39").Value = "3.57"
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("D41").Value = "14.38"
$ws.Range("E41").Value = "  -7.82%  "
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "1.743.73"
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("D45").Value = "82.82"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "99.97"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "4.93"
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "8.11"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "1.58"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "54.57"
$ws.Range("E51").Value = "  -2.70%  "
